# Generate Report for Handoff
#
# A new source file "ffff4a025b0b-2b31-4397-8070-53552b021016.md" has been
# marked "Ready for handoff". Also the previously-pending file's guid was
# regenerated from 680a7f29-605f-4fc4-917f-b40c26500715 to
# 6f416a5d-cbcb-44f0-a10f-8c2ecb1d8161 (and its handoff package hash / time
# refreshed). Each sheet (Overview, zh-cn, de-de) gets a new row inserted
# for the new file, pushing the ".localization-config" row down by one.

$wb = $excel.ActiveWorkbook

$oldGuid = "680a7f29-605f-4fc4-917f-b40c26500715"
$newGuid = "6f416a5d-cbcb-44f0-a10f-8c2ecb1d8161"
$newGuid2 = "ffff4a025b0b-2b31-4397-8070-53552b021016"

$oldHash = "f349857c19fd84875b3b5f6d4995882d680db096"
$newHash = "d1aff84ce4337b765bc52952f1d5edf069dfd7e8"

$zhFile = "$newGuid.$newHash.zh-cn.xlf"
$deFile = "$newGuid.$newHash.de-de.xlf"

$zhTime = "2016-03-10 09:43:00"
$deTime = "2016-03-10 09:43:10"

$epoch = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Hyperlinks.Delete()

$ws.Rows.Item(3).Insert()

$ws.Range("A2").Value = "$newGuid.md"
$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("C2").Value = "Ready for handoff"

$ws.Range("A3").Value = "$newGuid2.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"

$ws.Range("A4").Value = ".localization-config"
$ws.Range("B4").Value = "Not to be localized"
$ws.Range("C4").Value = "Not to be localized"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/6abf37440bebf913df272d613c6fd3331749a889/e2e/$newGuid.md", "", "", "$newGuid.md")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/6abf37440bebf913df272d613c6fd3331749a889/e2e/$newGuid2.md", "", "", "$newGuid2.md")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/ff75ff2df72d0bae9f84c41b053d24180944fde8/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Hyperlinks.Delete()

$ws.Rows.Item(3).Insert()

$ws.Range("A2").Value = "$newGuid.md"
$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("C2").Value = $zhFile
$ws.Range("D2").Value = $zhTime
$ws.Range("G2").Value = $epoch
$ws.Range("H2").Value = "Include"

$ws.Range("A3").Value = "$newGuid2.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = $zhFile
$ws.Range("D3").Value = $zhTime
$ws.Range("G3").Value = $epoch
$ws.Range("H3").Value = "Include"

$ws.Range("A4").Value = ".localization-config"
$ws.Range("B4").Value = "Not to be localized"
$ws.Range("D4").Value = $epoch
$ws.Range("G4").Value = $epoch
$ws.Range("H4").Value = "Ignored"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/6abf37440bebf913df272d613c6fd3331749a889/e2e/$newGuid.md", "", "", "$newGuid.md")
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/69edd87f0ac272b1962c64872ee39552c01518da/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhFile", "", "", $zhFile)
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/6abf37440bebf913df272d613c6fd3331749a889/e2e/$newGuid2.md", "", "", "$newGuid2.md")
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/69edd87f0ac272b1962c64872ee39552c01518da/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhFile", "", "", $zhFile)
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/ff75ff2df72d0bae9f84c41b053d24180944fde8/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Hyperlinks.Delete()

$ws.Rows.Item(3).Insert()

$ws.Range("A2").Value = "$newGuid.md"
$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("C2").Value = $deFile
$ws.Range("D2").Value = $deTime
$ws.Range("G2").Value = $epoch
$ws.Range("H2").Value = "Include"

$ws.Range("A3").Value = "$newGuid2.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = $deFile
$ws.Range("D3").Value = $deTime
$ws.Range("G3").Value = $epoch
$ws.Range("H3").Value = "Include"

$ws.Range("A4").Value = ".localization-config"
$ws.Range("B4").Value = "Not to be localized"
$ws.Range("D4").Value = $epoch
$ws.Range("G4").Value = $epoch
$ws.Range("H4").Value = "Ignored"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/6abf37440bebf913df272d613c6fd3331749a889/e2e/$newGuid.md", "", "", "$newGuid.md")
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0640527d5c2a6f87e8298bd9856860acdbfa8e6b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$deFile", "", "", $deFile)
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/6abf37440bebf913df272d613c6fd3331749a889/e2e/$newGuid2.md", "", "", "$newGuid2.md")
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0640527d5c2a6f87e8298bd9856860acdbfa8e6b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$deFile", "", "", $deFile)
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/ff75ff2df72d0bae9f84c41b053d24180944fde8/.localization-config", "", "", ".localization-config")

$wb.Worksheets.Item("Overview").Activate()
